# daily auto push: 2025-10-10 02:00 UTC
# Append the new daily log row (row 88) to the bottom of the sheet's
# existing data table: Date | Day-of-week | Time | Ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 88

# Column A holds a date-looking string ("2025/10/10") but, like every
# other row already in the sheet, it must stay plain text (not get
# auto-converted into an Excel date serial/number format). Leading the
# value with an apostrophe forces text entry; then resetting the style
# back to Normal strips the quote-prefix formatting Excel would
# otherwise stamp on the cell, so it ends up identical in style to its
# neighbours (no explicit style id).
$ws.Cells.Item($newRow, 1).Value = "'2025/10/10"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "金"
$ws.Cells.Item($newRow, 3).Value = 9
$ws.Cells.Item($newRow, 4).Value = 201
